# Auto-generated: apply market-data refresh values per diff (commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H70").Value = 1320.8077
$ws.Range("I70").Value = 1339.2084
$ws.Range("J70").Value = 1100
$ws.Range("K70").Value = 4017.6252
$ws.Range("L70").Value = 3300
$ws.Range("M70").Value = -3747.6252
$ws.Range("N70").Value = -3840
$ws.Range("H73").Value = 1320.8077
$ws.Range("I73").Value = 1339.2084
$ws.Range("J73").Value = 1100
$ws.Range("K73").Value = 4017.6252
$ws.Range("L73").Value = 3300
$ws.Range("M73").Value = -3081.6252
$ws.Range("N73").Value = -5172
$ws.Range("H86").Value = 185047.17
$ws.Range("I86").Value = 221260.6
$ws.Range("J86").Value = 3980
$ws.Range("K86").Value = 221260.6
$ws.Range("L86").Value = 3980
$ws.Range("M86").Value = -220137.6
$ws.Range("N86").Value = -6226
$ws.Range("H89").Value = 185047.17
$ws.Range("I89").Value = 221260.6
$ws.Range("J89").Value = 3980
$ws.Range("K89").Value = 1106303
$ws.Range("L89").Value = 19900
$ws.Range("M89").Value = -1100687
$ws.Range("N89").Value = -31132
$ws.Range("H98").Value = 2327.6
$ws.Range("I98").Value = 2557.5881
$ws.Range("J98").Value = 1838.875
$ws.Range("K98").Value = 2557.5881
$ws.Range("L98").Value = 1838.875
$ws.Range("M98").Value = -1059.5881
$ws.Range("N98").Value = -4834.875
$ws.Range("H122").Value = 2327.6
$ws.Range("I122").Value = 2557.5881
$ws.Range("J122").Value = 1838.875
$ws.Range("K122").Value = 7672.7643
$ws.Range("L122").Value = 5516.625
$ws.Range("M122").Value = -5222.7643
$ws.Range("N122").Value = -10416.625
$ws.Range("H137").Value = 1763.7858
$ws.Range("I137").Value = 2801.6875
$ws.Range("J137").Value = 1125.0769
$ws.Range("K137").Value = 8405.0625
$ws.Range("L137").Value = 3375.2307
$ws.Range("M137").Value = -5855.0625
$ws.Range("N137").Value = -8475.2307
$ws.Range("H138").Value = 2758.2097
$ws.Range("I138").Value = 1571.1143
$ws.Range("J138").Value = 4297.037
$ws.Range("K138").Value = 4713.3429
$ws.Range("L138").Value = 12891.111
$ws.Range("M138").Value = 426.6571000000004
$ws.Range("N138").Value = -23171.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4459.7476
$ws.Range("I32").Value = 4459.7476
$ws.Range("K32").Value = 4459.7476
$ws.Range("M32").Value = -4172.7476
$ws.Range("H74").Value = 2645.6843
$ws.Range("I74").Value = 3926.6667
$ws.Range("J74").Value = 2054.4614
$ws.Range("K74").Value = 3926.6667
$ws.Range("L74").Value = 2054.4614
$ws.Range("M74").Value = -3052.6667
$ws.Range("N74").Value = -3802.4614
$ws.Range("H77").Value = 2645.6843
$ws.Range("I77").Value = 3926.6667
$ws.Range("J77").Value = 2054.4614
$ws.Range("K77").Value = 19633.3335
$ws.Range("L77").Value = 10272.307
$ws.Range("M77").Value = -15265.3335
$ws.Range("N77").Value = -19008.307
$ws.Range("H97").Value = 1179
$ws.Range("I97").Value = 1042.1666
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1042.1666
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -546.1666
$ws.Range("N97").Value = -2992
$ws.Range("H102").Value = 1625
$ws.Range("I102").Value = 1192.7273
$ws.Range("K102").Value = 1192.7273
$ws.Range("M102").Value = 429.2727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H105").Value = 1876.7368
$ws.Range("I105").Value = 1792.5
$ws.Range("J105").Value = 2326
$ws.Range("K105").Value = 1792.5
$ws.Range("L105").Value = 2326
$ws.Range("M105").Value = -45.5
$ws.Range("N105").Value = -5820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 76924350
$ws.Range("I122").Value = 100000840
$ws.Range("J122").Value = 2733.3333
$ws.Range("K122").Value = 300002520
$ws.Range("L122").Value = 8199.999899999999
$ws.Range("M122").Value = -300000070
$ws.Range("N122").Value = -13099.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 774.9375
$ws.Range("I5").Value = 463.41666
$ws.Range("J5").Value = 1709.5
$ws.Range("K5").Value = 1390.24998
$ws.Range("L5").Value = 5128.5
$ws.Range("M5").Value = -1278.24998
$ws.Range("N5").Value = -5352.5
$ws.Range("H23").Value = 193.07692
$ws.Range("I23").Value = 143.16667
$ws.Range("J23").Value = 235.85715
$ws.Range("K23").Value = 429.50001
$ws.Range("L23").Value = 707.5714499999999
$ws.Range("M23").Value = -194.50001
$ws.Range("N23").Value = -1177.57145
$ws.Range("H135").Value = 774.9375
$ws.Range("I135").Value = 463.41666
$ws.Range("J135").Value = 1709.5
$ws.Range("K135").Value = 4170.74994
$ws.Range("L135").Value = 15385.5
$ws.Range("M135").Value = -1635.74994
$ws.Range("N135").Value = -20455.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 43686.555
$ws.Range("I140").Value = 31666.666
$ws.Range("J140").Value = 67726.336
$ws.Range("K140").Value = 31666.666
$ws.Range("L140").Value = 67726.336
$ws.Range("M140").Value = -26486.666
$ws.Range("N140").Value = -78086.336
$ws.Range("H141").Value = 67449.75
$ws.Range("J141").Value = 67449.75
$ws.Range("L141").Value = 67449.75
$ws.Range("N141").Value = -77809.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3982.8572
$ws.Range("I7").Value = 3500
$ws.Range("J7").Value = 4626.6665
$ws.Range("K7").Value = 3500
$ws.Range("L7").Value = 4626.6665
$ws.Range("M7").Value = -3388
$ws.Range("N7").Value = -4850.6665
$ws.Range("H40").Value = 125005750
$ws.Range("I40").Value = 250004260
$ws.Range("J40").Value = 7246.25
$ws.Range("K40").Value = 250004260
$ws.Range("L40").Value = 7246.25
$ws.Range("M40").Value = -250004124
$ws.Range("N40").Value = -7518.25
$ws.Range("H55").Value = 255.93103
$ws.Range("I55").Value = 276.9375
$ws.Range("J55").Value = 230.07692
$ws.Range("K55").Value = 276.9375
$ws.Range("L55").Value = 230.07692
$ws.Range("M55").Value = -103.9375
$ws.Range("N55").Value = -576.07692
$ws.Range("H126").Value = 3982.8572
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 4626.6665
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 13879.9995
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -18819.9995
$ws.Range("H136").Value = 1750.0667
$ws.Range("I136").Value = 1476.4828
$ws.Range("J136").Value = 2245.9375
$ws.Range("K136").Value = 4429.4484
$ws.Range("L136").Value = 6737.8125
$ws.Range("M136").Value = -1879.4484
$ws.Range("N136").Value = -11837.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 4000
$ws.Range("J5").Value = 4000
$ws.Range("L5").Value = 4000
$ws.Range("N5").Value = -4224
$ws.Range("H62").Value = 3808.3333
$ws.Range("I62").Value = 2745.5
$ws.Range("J62").Value = 5934
$ws.Range("K62").Value = 2745.5
$ws.Range("L62").Value = 5934
$ws.Range("M62").Value = -2121.5
$ws.Range("N62").Value = -7182
$ws.Range("H65").Value = 3808.3333
$ws.Range("I65").Value = 2745.5
$ws.Range("J65").Value = 5934
$ws.Range("K65").Value = 13727.5
$ws.Range("L65").Value = 29670
$ws.Range("M65").Value = -10607.5
$ws.Range("N65").Value = -35910
$ws.Range("H122").Value = 4340
$ws.Range("I122").Value = 6000
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -15550
$ws.Range("N122").Value = -10450
$ws.Range("H132").Value = 2335.32
$ws.Range("I132").Value = 1961.375
$ws.Range("J132").Value = 3000.111
$ws.Range("K132").Value = 5884.125
$ws.Range("L132").Value = 9000.332999999999
$ws.Range("M132").Value = -3354.125
$ws.Range("N132").Value = -14060.333
$ws.Range("H140").Value = 37885.8
$ws.Range("J140").Value = 37885.8
$ws.Range("L140").Value = 37885.8
$ws.Range("N140").Value = -48245.8
$ws.Range("H141").Value = 41500
$ws.Range("J141").Value = 41500
$ws.Range("L141").Value = 41500
$ws.Range("N141").Value = -51860
